# Revert "Merging 0.1.8 w VitalSigns"
#
# This reverts the Metadata sheet's Version/Status/Date/Contact values back
# to their pre-merge values, removes the (now redundant) "Jurisdiction" row,
# and renames the second worksheet back to "Include from LOINC".

$wb = $excel.ActiveWorkbook

# --- Rename the "Include #0" sheet back to "Include from LOINC" ---
$wsInclude = $wb.Worksheets.Item(2)
$wsInclude.Name = "Include from LOINC"

# --- Update the Metadata sheet ---
$ws = $wb.Worksheets.Item(1)

# Version: 0.1.8 -> 0.1.6
$ws.Range("B3").Value = "0.1.6"

# Status: draft -> active
$ws.Range("B6").Value = "active"

# Date: 2025-02-18T16:25:42-06:00 -> 2023-05-05T10:50:04-05:00
$ws.Range("B8").Value = "2023-05-05T10:50:04-05:00"

# Contact rows: both collapse to "No display for ContactDetail"
$ws.Range("B10").Value = "No display for ContactDetail"
$ws.Range("B11").Value = "No display for ContactDetail"

# Remove the "Jurisdiction" row entirely (row 12), shifting later rows up
$ws.Rows.Item(12).Delete()
